# CIDC-1410 add fix for clinical data participant counting
#
# The "cimac_part_id" column (and the participant-id values under it) was
# the FIRST column (A) in the clinical data sheet, which was throwing off
# downstream participant counting. This moves it to be the LAST column
# (G) instead, shifting every other column one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Snapshot the current 5x7 block (header row + 4 data rows) before
# overwriting anything, then rewrite each row shifted one column to the
# left, wrapping the old column A into the new column G.
$lastCol = 7
$lastRow = 5

for ($r = 1; $r -le $lastRow; $r++) {
    $rowValues = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowValues += $ws.Cells.Item($r, $c).Value2
    }

    # Rotate left: new col c gets old col (c+1), new last col gets old col 1
    for ($c = 1; $c -le $lastCol; $c++) {
        $srcIndex = $c % $lastCol
        $ws.Cells.Item($r, $c).Value = $rowValues[$srcIndex]
    }
}

# Restore the header row's formatting/height (Excel re-measures this to
# 27.95 after the edit in the source file).
$ws.Rows.Item(1).RowHeight = 27.95

# Match the saved selection: column G (the relocated cimac_part_id column)
# for the data rows.
$ws.Range("G2:G5").Select() | Out-Null
